# Reorder the player data rows (A2:C16) to match the updated sheet.
# Row 1 (header) and row 17 (CJ McCollum) stay unchanged; rows 2-16 are
# re-populated with the same set of player records but in a new order.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("Dennis Schröder", "PG", "Brooklyn Nets"),
    @("Kyrie Irving", "PG,SG", "Dallas Mavericks"),
    @("Jordan Poole", "PG,SG", "Washington Wizards"),
    @("Brandon Boston Jr.", "SG,SF", "New Orleans Pelicans"),
    @("Christian Braun", "SG,SF", "Denver Nuggets"),
    @("Zach LaVine", "SG,SF", "Chicago Bulls"),
    @("Joel Embiid", "C", "Philadelphia 76ers"),
    @("Lauri Markkanen", "SF,PF", "Utah Jazz"),
    @("Keyonte George", "PG,SG", "Utah Jazz"),
    @("John Collins", "PF,C", "Utah Jazz"),
    @("Tobias Harris", "SF,PF", "Detroit Pistons"),
    @("Shai Gilgeous-Alexander", "PG", "Oklahoma City Thunder"),
    @("RJ Barrett", "SF,PF", "Toronto Raptors"),
    @("Jimmy Butler", "SF,PF", "Miami Heat"),
    @("Jalen Williams", "SG,SF,PF,C", "Oklahoma City Thunder")
)

$row = 2
foreach ($record in $data) {
    $ws.Cells.Item($row, 1).Value = $record[0]
    $ws.Cells.Item($row, 2).Value = $record[1]
    $ws.Cells.Item($row, 3).Value = $record[2]
    $row++
}
